# Auto-generated edit script applying the diff to Sheets/Cuchulainn_Profits.xlsx
# (workbook tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 1182.5
$ws.Range("I115").Value = 1182.5
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 3547.5
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1980.5
$ws.Range("N115").ClearContents()
$ws.Range("H125").Value = 1500
$ws.Range("J125").Value = 1500
$ws.Range("L125").Value = 13500
$ws.Range("N125").Value = -18420
$ws.Range("H127").Value = 4999.6665
$ws.Range("I127").Value = 4999.6665
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 14998.9995
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -10038.9995
$ws.Range("N127").ClearContents()
$ws.Range("H137").Value = 9447.182000000001
$ws.Range("I137").Value = 6864.875
$ws.Range("K137").Value = 20594.625
$ws.Range("M137").Value = -18044.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3512.2222
$ws.Range("I2").Value = 1870
$ws.Range("K2").Value = 1870
$ws.Range("M2").Value = -1757
$ws.Range("H61").Value = 5138.615
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 5714.2856
$ws.Range("J74").Value = 7400
$ws.Range("L74").Value = 7400
$ws.Range("N74").Value = -9148
$ws.Range("H77").Value = 5714.2856
$ws.Range("J77").Value = 7400
$ws.Range("L77").Value = 37000
$ws.Range("N77").Value = -45736
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H110").Value = 2000
$ws.Range("I110").Value = 2000
$ws.Range("K110").Value = 2000
$ws.Range("M110").Value = 45
$ws.Range("H116").Value = 3512.2222
$ws.Range("I116").Value = 1870
$ws.Range("K116").Value = 1870
$ws.Range("M116").Value = 424
$ws.Range("H132").Value = 6559.222
$ws.Range("I132").Value = 996.6
$ws.Range("J132").Value = 13512.5
$ws.Range("K132").Value = 2989.8
$ws.Range("L132").Value = 40537.5
$ws.Range("M132").Value = -459.8000000000002
$ws.Range("N132").Value = -45597.5
$ws.Range("H136").Value = 5138.615
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3512.2222
$ws.Range("I3").Value = 1870
$ws.Range("K3").Value = 1870
$ws.Range("M3").Value = -1756
$ws.Range("H86").Value = 775
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 775
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 636.8333
$ws.Range("I107").Value = 622.2
$ws.Range("K107").Value = 622.2
$ws.Range("M107").Value = 1297.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4566.0386
$ws.Range("I31").Value = 3512.2632
$ws.Range("K31").Value = 3512.2632
$ws.Range("M31").Value = -3217.2632
$ws.Range("H34").Value = 4566.0386
$ws.Range("I34").Value = 3512.2632
$ws.Range("K34").Value = 3512.2632
$ws.Range("M34").Value = -3310.2632
$ws.Range("H58").Value = 4820.6665
$ws.Range("I58").Value = 1231
$ws.Range("J58").Value = 12000
$ws.Range("K58").Value = 1231
$ws.Range("L58").Value = 12000
$ws.Range("M58").Value = -1028
$ws.Range("N58").Value = -12406
$ws.Range("H106").Value = 25223.666
$ws.Range("J106").Value = 25223.666
$ws.Range("L106").Value = 25223.666
$ws.Range("N106").Value = -27747.666
$ws.Range("H134").Value = 3239.9375
$ws.Range("J134").Value = 10006.75
$ws.Range("L134").Value = 30020.25
$ws.Range("N134").Value = -35090.25
$ws.Range("H136").Value = 4820.6665
$ws.Range("I136").Value = 1231
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 3693
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -1143
$ws.Range("N136").Value = -41100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 166.66667
$ws.Range("I7").Value = 149.5
$ws.Range("J7").Value = 201
$ws.Range("K7").Value = 448.5
$ws.Range("L7").Value = 603
$ws.Range("M7").Value = -336.5
$ws.Range("N7").Value = -827
$ws.Range("H11").Value = 483.33334
$ws.Range("J11").Value = 1000
$ws.Range("L11").Value = 3000
$ws.Range("N11").Value = -3280
$ws.Range("H13").Value = 1479.8
$ws.Range("I13").Value = 1549.6666
$ws.Range("J13").Value = 1375
$ws.Range("K13").Value = 4648.9998
$ws.Range("L13").Value = 4125
$ws.Range("M13").Value = -4480.9998
$ws.Range("N13").Value = -4461
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 137
$ws.Range("N15").ClearContents()
$ws.Range("H34").Value = 1250
$ws.Range("I34").Value = 83.333336
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 250.000008
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -166.000008
$ws.Range("N34").Value = -9168
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H108").Value = 165.75
$ws.Range("I108").Value = 165.75
$ws.Range("K108").Value = 497.25
$ws.Range("M108").Value = 2382.75
$ws.Range("H114").Value = 915.25
$ws.Range("J114").Value = 1120.3334
$ws.Range("L114").Value = 3361.0002
$ws.Range("N114").Value = -9869.0002
$ws.Range("H117").Value = 2202.2222
$ws.Range("J117").Value = 2677.3333
$ws.Range("L117").Value = 8031.999899999999
$ws.Range("N117").Value = -14915.9999
$ws.Range("H122").Value = 1000.25
$ws.Range("J122").Value = 1074.5
$ws.Range("L122").Value = 9670.5
$ws.Range("N122").Value = -14570.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5174.4375
$ws.Range("J132").Value = 12198
$ws.Range("L132").Value = 36594
$ws.Range("N132").Value = -41654

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1249.875
$ws.Range("I55").Value = 1166.5
$ws.Range("K55").Value = 1166.5
$ws.Range("M55").Value = -993.5
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H132").Value = 13388.833
$ws.Range("I132").Value = 8953.857
$ws.Range("J132").Value = 19597.8
$ws.Range("K132").Value = 26861.571
$ws.Range("L132").Value = 58793.39999999999
$ws.Range("M132").Value = -24331.571
$ws.Range("N132").Value = -63853.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11286.077
$ws.Range("I132").Value = 10090.25
$ws.Range("K132").Value = 30270.75
$ws.Range("M132").Value = -27740.75
$ws.Range("H136").Value = 10232.833
$ws.Range("I136").Value = 5685
$ws.Range("K136").Value = 17055
$ws.Range("M136").Value = -14505

Write-Output "Applied all cell updates."